# Bug fixes when adding a point
# Two new records were inserted at the top of the data table (right under
# the header row), pushing the existing rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right after the header row (row 1), so the
# new points land on rows 2 and 3 and all prior data shifts down by 2.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# The inserted rows pick up formatting from the row above (the bold /
# bordered header style); clear that so they match the other plain data
# rows.
$ws.Rows.Item(2).ClearFormats()
$ws.Rows.Item(3).ClearFormats()

# New row 2: BEAML-RH point
$ws.Range("A2").Value = "BEAML-RH "
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 6475
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "Right"

# New row 3: BEAMU-RH point
$ws.Range("A3").Value = "BEAMU-RH "
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 4961
$ws.Range("D3").Value = "Delaminacion"
$ws.Range("E3").Value = "Right"
